# Workplace Ergonomics doc update
#  1. Merge the title's split runs ("<New Lands" / ">  -" / " " / "Workplace
#     Ergonomics") into a single run of text, removing the stray gramStart/
#     gramEnd proofing marks in the process.
#  2. Add a new bold "BitBucket" paragraph straight after the GitHub
#     screenshot (and before the "Operating System Settings" heading), to
#     match the new Trello / GitHub / BitBucket trio of third-party tools.

$d = $word.ActiveDocument

# --- 1. Title fix -----------------------------------------------------
$find = $d.Content.Find
$find.Execute(
    "<New Lands>  - Workplace Ergonomics",  # FindText
    $true,                                  # MatchCase
    $false,                                 # MatchWholeWord
    $false,                                 # MatchWildcards
    $false,                                 # MatchSoundsLike
    $false,                                 # MatchAllWordForms
    $true,                                  # Forward
    1,                                      # Wrap (wdFindContinue)
    $false,                                 # Format
    "<New Lands>  - Workplace Ergonomics",  # ReplaceWith (same text - collapses the runs)
    2                                       # Replace (wdReplaceAll)
) | Out-Null

# --- 2. Add the BitBucket paragraph ------------------------------------
# Find the paragraph that holds the GitHub picture: it's the last empty
# paragraph (only a drawing, no text) before "Operating System Settings".
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13) -eq "") {
        if (($i + 1) -le $d.Paragraphs.Count) {
            $nextPara = $d.Paragraphs.Item($i + 1)
            if ($nextPara.Range.Text.TrimEnd([char]13) -eq "Operating System Settings") {
                $targetIndex = $i
                break
            }
        }
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the GitHub picture paragraph"
}

$githubPicPara = $d.Paragraphs.Item($targetIndex)
$githubPicPara.Range.InsertParagraphAfter() | Out-Null

$bitBucketPara = $d.Paragraphs.Item($targetIndex + 1)
$bitBucketPara.Range.Text = "BitBucket"

Write-Output "done"
